$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "91.737.96"
$ws.Range("E2").Value = "  +1.27%  "
$ws.Range("D3").Value = "3.123.83"
$ws.Range("E3").Value = "  +1.78%  "
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").Value = "246.38"
$ws.Range("E5").Value = "  +1.42%  "
$ws.Range("D6").Value = "618.02"
$ws.Range("E6").Value = "  +0.36%  "
$ws.Range("E7").Value = "  -1.77%  "
$ws.Range("D8").Value = "0.384"
$ws.Range("E8").Value = "  +5.53%  "
$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  -0.15%  "
$ws.Range("D10").Value = "3.120.19"
$ws.Range("E10").Value = "  +1.64%  "
$ws.Range("D11").Value = "0.736"
$ws.Range("E11").Value = "  -0.15%  "
$ws.Range("E12").Value = "  +1.76%  "
$ws.Range("D13").Value = "0.0000251"
$ws.Range("E13").Value = "  +1.67%  "
$ws.Range("D14").Value = "5.60"
$ws.Range("E14").Value = "  +3.29%  "
$ws.Range("D15").Value = "34.83"
$ws.Range("E15").Value = "  +0.08%  "
$ws.Range("D16").Value = "91.566.38"
$ws.Range("E16").Value = "  +0.68%  "
$ws.Range("D17").Value = "3.703.21"
$ws.Range("E17").Value = "  +1.29%  "
$ws.Range("D18").Value = "3.071.63"
$ws.Range("E18").Value = "  -0.83%  "
$ws.Range("D19").Value = "3.71"
$ws.Range("E19").Value = "  +0.99%  "
$ws.Range("D20").Value = "14.92"
$ws.Range("E20").Value = "  +3.68%  "
$ws.Range("B21").Value = "Polkadot"
$ws.Range("C21").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D21").Value = "5.84"
$ws.Range("E21").Value = "  +1.83%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "9.55"
$ws.Range("E22").Value = "  +6.04%  "
$ws.Range("D23").Value = "447.88"
$ws.Range("E23").Value = "  +2.00%  "
$ws.Range("E24").Value = "  -3.90%  "
$ws.Range("D25").Value = "5.88"
$ws.Range("E25").Value = "  +5.51%  "
$ws.Range("D26").Value = "88.01"
$ws.Range("E26").Value = "  -3.02%  "
$ws.Range("D27").Value = "11.78"
$ws.Range("E27").Value = "  +0.41%  "
$ws.Range("D29").Value = "0.146"
$ws.Range("E29").Value = "  +32.15%  "
$ws.Range("E30").Value = "  +0.20%  "
$ws.Range("D31").Value = "0.237"
$ws.Range("E31").Value = "  -5.55%  "
$ws.Range("E32").Value = "  -7.81%  "
$ws.Range("E33").Value = "  +5.93%  "
$ws.Range("D34").Value = "9.34"
$ws.Range("E34").Value = "  +2.77%  "
$ws.Range("E35").Value = "  -1.01%  "
$ws.Range("D36").Value = "7.90"
$ws.Range("E36").Value = "  +2.93%  "
$ws.Range("D37").Value = "26.29"
$ws.Range("E37").Value = "  +0.16%  "
$ws.Range("D38").Value = "4.22"
$ws.Range("E38").Value = "  +1.03%  "
$ws.Range("E39").Value = "  +1.82%  "
$ws.Range("D40").Value = "491.77"
$ws.Range("E40").Value = "  +0.47%  "
$ws.Range("D41").Value = "1.31"
$ws.Range("E41").Value = "  +1.92%  "
$ws.Range("D42").Value = "0.441"
$ws.Range("E42").Value = "  +6.44%  "
$ws.Range("E43").Value = "  -5.24%  "
$ws.Range("D44").Value = "22.18"
$ws.Range("E44").Value = "  +0.28%  "
$ws.Range("E45").Value = "  -0.05%  "
$ws.Range("D46").Value = "157.88"
$ws.Range("E46").Value = "  +2.57%  "
$ws.Range("D47").Value = "0.708"
$ws.Range("E47").Value = "  +4.07%  "
$ws.Range("E48").Value = "  +2.13%  "
$ws.Range("D49").Value = "1.36"
$ws.Range("E49").Value = "  +2.92%  "
$ws.Range("B50").Value = "Filecoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D50").Value = "4.42"
$ws.Range("E50").Value = "  -0.10%  "
$ws.Range("B51").Value = "OKB"
$ws.Range("C51").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D51").Value = "44.04"
$ws.Range("E51").Value = "  +0.04%  "
